$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the XPath values: drop the "[1]" index from each Table[1] reference
$ws.Range("B2").Value = "/NewDataSet/Table/Town"
$ws.Range("B3").Value = "/NewDataSet/Table/County"
$ws.Range("B4").Value = "/NewDataSet/Table/PostCode"

# Adjust column widths to match the new best-fit sizes
# (engine quantizes ColumnWidth to 1/6-character steps, so feed the
# pre-quantized values that land closest to the target stored widths)
$ws.Columns.Item(2).ColumnWidth = 27.666666666666668
$ws.Columns.Item(3).ColumnWidth = 12.5

# Update the active selection on the sheet
$ws.Range("E5").Select()
